$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: updated hydrogen demand value
$ws.Range("B3").Value = 15966399.99901282

# D3: value removed (cell becomes empty / inline string placeholder)
$ws.Range("D3").Value = ""

# C4: updated methanol demand value
$ws.Range("C4").Value = 5274.663204186195

# C5: updated ammonia demand value
$ws.Range("C5").Value = 13954.48126513115

# Row 7 label renamed from "Other" to "Biogas", value updated
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 2029.410957721489

# New row 8: "Other" category (moved down from row 7), with its own value
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats - match row label styling (style 1)
$excel.CutCopyMode = $false
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 177.4265722087592
